$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.906.08"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "3.738.95"
$ws.Range("E3").Value = "  -2.05%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.90"
$ws.Range("E5").Value = "  -1.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.85"
$ws.Range("E6").Value = "  -1.64%  "
$ws.Range("D7").Value = "3.741.03"
$ws.Range("E7").Value = "  -1.98%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  -2.03%  "
$ws.Range("E10").Value = "  -3.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.51"
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("E12").Value = "  -2.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000263"
$ws.Range("E13").Value = "  -4.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.50"
$ws.Range("E14").Value = "  -1.56%  "
$ws.Range("D15").Value = "4.363.38"
$ws.Range("E15").Value = "  -2.11%  "
$ws.Range("D16").Value = "3.725.37"
$ws.Range("E16").Value = "  -2.18%  "
$ws.Range("D17").Value = "67.894.61"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.29"
$ws.Range("E18").Value = "  -2.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.06"
$ws.Range("E19").Value = "  -5.60%  "
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.79"
$ws.Range("E21").Value = "  -0.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "468.27"
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.703"
$ws.Range("E23").Value = "  -5.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.91"
$ws.Range("E24").Value = "  -0.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.23"
$ws.Range("E25").Value = "  -3.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000136"
$ws.Range("E26").Value = "  -11.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.07"
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.18"
$ws.Range("E28").Value = "  -1.83%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "3.881.77"
$ws.Range("E30").Value = "  -2.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.78"
$ws.Range("E31").Value = "  -4.97%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.35"
$ws.Range("E32").Value = "  -5.25%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.24"
$ws.Range("E33").Value = "  -1.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.81"
$ws.Range("E34").Value = "  -3.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E36").Value = "  -3.13%  "
$ws.Range("D37").Value = "3.689.17"
$ws.Range("E37").Value = "  -2.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.102"
$ws.Range("E38").Value = "  -4.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.44"
$ws.Range("E39").Value = "  -10.61%  "
$ws.Range("E40").Value = "  -0.84%  "
$ws.Range("E41").Value = "  -1.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.77"
$ws.Range("E42").Value = "  -3.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.307"
$ws.Range("E45").Value = "  -3.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.59"
$ws.Range("E46").Value = "  -2.57%  "
$ws.Range("E47").Value = "  -2.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "45.40"
$ws.Range("E48").Value = "  -2.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "391.97"
$ws.Range("E49").Value = "  -4.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "143.56"
$ws.Range("E50").Value = "  +0.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "25.53"
$ws.Range("E51").Value = "  +0.29%  "
